# Updated cryptos list on Fri Jul  5 06:12:32 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) snapshot values
# pulled from coinranking.com, and reorders two rows (40/41) whose ranking
# swapped: "RenzoRestakedETH" now ranks above "Hedera".
#
# Price/Volume cells are stored as literal text (e.g. "54.565.92",
# "  -7.27%  "), not numbers, so assignment uses Range.Value with the
# exact display string. Some new Price strings (e.g. "1.00", "477.78")
# parse as plain numbers; a leading apostrophe forces Excel to keep them
# as text (matching the source workbook), and the quote-prefix cell
# styling that introduces is cleared afterwards via Style = "Normal" so
# the cells keep their original (default) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range('D2').Value = '54.565.92'
$ws.Range('E2').Value = '  -7.27%  '

# Row 3
$ws.Range('D3').Value = '2.878.17'
$ws.Range('E3').Value = '  -10.60%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').Value = '''477.78'
$ws.Range('E5').Value = '  -11.06%  '

# Row 6
$ws.Range('D6').Value = '''126.31'
$ws.Range('E6').Value = '  -6.76%  '

# Row 7
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  -0.05%  '

# Row 8
$ws.Range('D8').Value = '2.873.02'
$ws.Range('E8').Value = '  -10.75%  '

# Row 9
$ws.Range('D9').Value = '''0.404'
$ws.Range('E9').Value = '  -11.97%  '

# Row 10
$ws.Range('D10').Value = '''6.68'
$ws.Range('E10').Value = '  -11.75%  '

# Row 11
$ws.Range('D11').Value = '''0.0969'
$ws.Range('E11').Value = '  -15.31%  '

# Row 12
$ws.Range('D12').Value = '''0.331'
$ws.Range('E12').Value = '  -15.83%  '

# Row 13
$ws.Range('E13').Value = '  -4.22%  '

# Row 14
$ws.Range('D14').Value = '3.376.43'
$ws.Range('E14').Value = '  -10.62%  '

# Row 15
$ws.Range('D15').Value = '''22.83'
$ws.Range('E15').Value = '  -11.71%  '

# Row 16
$ws.Range('D16').Value = '54.513.18'
$ws.Range('E16').Value = '  -7.48%  '

# Row 17
$ws.Range('D17').Value = '2.889.09'
$ws.Range('E17').Value = '  -10.29%  '

# Row 18
$ws.Range('D18').Value = '''0.0000135'
$ws.Range('E18').Value = '  -14.54%  '

# Row 19
$ws.Range('D19').Value = '''5.25'
$ws.Range('E19').Value = '  -11.33%  '

# Row 20
$ws.Range('D20').Value = '''11.53'
$ws.Range('E20').Value = '  -13.21%  '

# Row 21
$ws.Range('D21').Value = '''7.07'
$ws.Range('E21').Value = '  -14.02%  '

# Row 22
$ws.Range('D22').Value = '''305.97'
$ws.Range('E22').Value = '  -15.65%  '

# Row 23
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.11%  '

# Row 24
$ws.Range('D24').Value = '''0.446'
$ws.Range('E24').Value = '  -14.18%  '

# Row 25
$ws.Range('D25').Value = '''59.33'
$ws.Range('E25').Value = '  -15.79%  '

# Row 26
$ws.Range('D26').Value = '''0.998'
$ws.Range('E26').Value = '  -0.18%  '

# Row 27
$ws.Range('D27').Value = '''0.153'
$ws.Range('E27').Value = '  -10.45%  '

# Row 28
$ws.Range('D28').Value = '''0.997'
$ws.Range('E28').Value = '  -0.24%  '

# Row 29
$ws.Range('D29').Value = '0.0₃0817'
$ws.Range('E29').Value = '  -14.98%  '

# Row 30
$ws.Range('D30').Value = '''6.30'
$ws.Range('E30').Value = '  -11.46%  '

# Row 31
$ws.Range('D31').Value = '''1.13'
$ws.Range('E31').Value = '  -6.37%  '

# Row 32
$ws.Range('D32').Value = '''6.19'
$ws.Range('E32').Value = '  -12.39%  '

# Row 33
$ws.Range('D33').Value = '''19.07'
$ws.Range('E33').Value = '  -12.71%  '

# Row 34
$ws.Range('E34').Value = '  -16.31%  '

# Row 35
$ws.Range('D35').Value = '''4.22'
$ws.Range('E35').Value = '  -14.26%  '

# Row 36
$ws.Range('D36').Value = '''136.30'
$ws.Range('E36').Value = '  -15.38%  '

# Row 37
$ws.Range('D37').Value = '''5.42'
$ws.Range('E37').Value = '  -15.42%  '

# Row 38
$ws.Range('D38').Value = '''1.21'
$ws.Range('E38').Value = '  -16.19%  '

# Row 39
$ws.Range('D39').Value = '''22.81'
$ws.Range('E39').Value = '  -12.73%  '

# Row 40
$ws.Range('B40').Value = 'RenzoRestakedETH'
$ws.Range('C40').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D40').Value = '2.910.04'
$ws.Range('E40').Value = '  -10.50%  '

# Row 41
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '''0.0616'
$ws.Range('E41').Value = '  -12.81%  '

# Row 42
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -0.07%  '

# Row 43
$ws.Range('D43').Value = '''35.61'
$ws.Range('E43').Value = '  -13.37%  '

# Row 44
$ws.Range('D44').Value = '''0.968'
$ws.Range('E44').Value = '  -12.30%  '

# Row 45
$ws.Range('D45').Value = '''0.601'
$ws.Range('E45').Value = '  -15.96%  '

# Row 46
$ws.Range('D46').Value = '''1.32'
$ws.Range('E46').Value = '  -11.99%  '

# Row 47
$ws.Range('D47').Value = '''3.41'
$ws.Range('E47').Value = '  -15.13%  '

# Row 48
$ws.Range('D48').Value = '2.049.42'
$ws.Range('E48').Value = '  -10.97%  '

# Row 49
$ws.Range('D49').Value = '''5.30'
$ws.Range('E49').Value = '  -15.86%  '

# Row 50
$ws.Range('D50').Value = '''17.90'
$ws.Range('E50').Value = '  -14.16%  '

# Row 51
$ws.Range('D51').Value = '''0.0212'
$ws.Range('E51').Value = '  -11.93%  '

# Clear the quote-prefix formatting the leading apostrophes above introduced,
# restoring the default (unstyled) cell format used by the rest of the sheet.
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
